$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 441, pushing existing rows 441-467 down to 442-468
$ws.Rows.Item(441).Insert()

# Populate the newly inserted row 441 with the new record
$ws.Cells.Item(441, 1).Value2  = 4
$ws.Cells.Item(441, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(441, 3).Value2  = "Los Lagos"
$ws.Cells.Item(441, 4).Value2  = 45267
$ws.Cells.Item(441, 5).Value2  = 10
$ws.Cells.Item(441, 6).Value2  = 100112028
$ws.Cells.Item(441, 7).Value2  = "Sandia"
$ws.Cells.Item(441, 8).Value2  = "Sin especificar"
$ws.Cells.Item(441, 9).Value2  = "Primera"
$ws.Cells.Item(441, 10).Value2 = 300
$ws.Cells.Item(441, 11).Value2 = 1100
$ws.Cells.Item(441, 12).Value2 = 1100
$ws.Cells.Item(441, 13).Value2 = 1100
$ws.Cells.Item(441, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(441, 15).Value2 = "Perú"
$ws.Cells.Item(441, 16).Value2 = 1100
$ws.Cells.Item(441, 17).Value2 = 1
$ws.Cells.Item(441, 18).Value2 = "Hortaliza"
